#
# Refresh the "cryptos" price/volume snapshot with the latest values from the feed
# (GitHub Actions scheduled update). Coin names/links are unchanged except that
# Algorand and TheSandbox swapped table positions (rows 40/41).
#
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain (non-ambiguous) value updates: Coin/Link text swaps and Volume(1h) percentages,
# plus Price values that already contain two dots so Excel keeps them as text automatically ---
$plainUpdates = @{
    "D2" = '26.821.32'
    "E2" = '  -2.50%  '
    "D3" = '1.775.45'
    "E3" = '  -3.09%  '
    "E4" = '  +0.60%  '
    "E5" = '  +0.57%  '
    "E6" = '  -1.69%  '
    "E7" = '  -1.54%  '
    "E8" = '  -1.35%  '
    "E9" = '  -1.83%  '
    "E10" = '  -3.70%  '
    "E11" = '  -2.30%  '
    "D12" = '1.807.02'
    "E12" = '  -0.52%  '
    "E13" = '  -3.00%  '
    "E14" = '  -3.01%  '
    "E15" = '  -1.96%  '
    "E16" = '  +0.72%  '
    "E17" = '  -1.76%  '
    "E18" = '  -2.92%  '
    "E19" = '  +0.44%  '
    "E20" = '  -3.29%  '
    "D21" = '26.727.60'
    "E21" = '  -2.47%  '
    "E22" = '  -2.58%  '
    "E23" = '  +1.53%  '
    "D24" = '1.973.97'
    "E24" = '  -3.95%  '
    "E25" = '  -3.15%  '
    "E26" = '  -0.77%  '
    "E27" = '  -4.80%  '
    "E28" = '  -2.29%  '
    "E29" = '  -0.16%  '
    "E30" = '  -11.91%  '
    "E31" = '  +0.75%  '
    "E32" = '  -4.78%  '
    "E33" = '  -4.51%  '
    "E34" = '  -5.03%  '
    "E35" = '  -4.44%  '
    "E36" = '  +0.56%  '
    "E37" = '  -1.30%  '
    "E38" = '  -2.64%  '
    "E39" = '  -4.81%  '
    "B40" = 'TheSandbox'
    "C40" = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
    "E40" = '  -3.80%  '
    "B41" = 'Algorand'
    "C41" = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    "E41" = '  -3.73%  '
    "E42" = '  -9.98%  '
    "E43" = '  -7.45%  '
    "E44" = '  -5.54%  '
    "E45" = '  +0.65%  '
    "E47" = '  -4.21%  '
    "E48" = '  -4.30%  '
    "E49" = '  -4.69%  '
    "E50" = '  -3.10%  '
    "E51" = '  -1.86%  '
}
foreach ($cellRef in $plainUpdates.Keys) {
    $ws.Range($cellRef).Value = $plainUpdates[$cellRef]
}

# --- Price values that look like plain numbers to Excel (single decimal point). ---
# Force the cell to Text format before assignment so the digits/format from the source feed
# are preserved exactly (e.g. '307.48' instead of being parsed into the number 307.48),
# then restore the cell's original (Normal) style so no formatting changes are introduced.
$textPriceUpdates = @{
    "D6" = '307.48'
    "D7" = '0.4223'
    "D8" = '0.3614'
    "D9" = '0.07163'
    "D10" = '0.8352'
    "D11" = '20.23'
    "D13" = '5.251'
    "D14" = '6.334'
    "D15" = '0.06803'
    "D18" = '0.000008671'
    "D20" = '14.94'
    "D22" = '5.008'
    "D23" = '11.03'
    "D25" = '1.919'
    "D26" = '153.35'
    "D27" = '18.09'
    "D28" = '5.029'
    "D29" = '114.30'
    "D30" = '1.621'
    "D31" = '0.08943'
    "D32" = '0.7188'
    "D33" = '2.847'
    "D34" = '4.320'
    "D35" = '1.089'
    "D37" = '1.077'
    "D38" = '0.01890'
    "D39" = '0.05077'
    "D40" = '0.4910'
    "D41" = '0.1605'
    "D42" = '2.521'
    "D43" = '6.095'
    "D44" = '7.911'
    "D45" = '1.007'
    "D46" = '104.61'
    "D47" = '10.06'
    "D48" = '0.06223'
    "D49" = '0.4472'
    "D50" = '1.569'
    "D51" = '1.711'
}
foreach ($cellRef in $textPriceUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $textPriceUpdates[$cellRef]
    $cell.Style = "Normal"
}
